{"js": "// Commit: \"Added cloud computing exp\"\n// 1) Add a new bulleted \"IBM DB2 on Cloud: cloud computing \" line to the\n//    \"Tools:\" list (right after the \"Tableau: data visualizations\" bullet,\n//    sharing the same numbering list / style as its sibling bullets).\n// 2) Remove the (now redundant) first blank spacer paragraph that sat right\n//    after the \"There are two different python scripts...\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// --- Part 1: insert the new \"IBM DB2 on Cloud\" bullet -------------------\nconst tableauPara = paragraphs.items.find(\n  (p) => p.style === \"List Paragraph\" && p.text.trim() === \"Tableau: data visualizations\"\n);\nif (!tableauPara) {\n  throw new Error(\"Could not find the 'Tableau: data visualizations' bullet paragraph\");\n}\n\n// Read the numbering list this paragraph belongs to so the new bullet joins\n// the very same list (numId 11 in the underlying OOXML).\nconst tableauList = tableauPara.list;\ntableauList.load(\"id\");\nawait context.sync();\nconst listId = tableauList.id;\n\nconst newBullet = tableauPara.insertParagraph(\"IBM DB2 on Cloud: cloud computing \", \"After\");\nnewBullet.styleBuiltIn = Word.BuiltInStyleName.listParagraph;\nnewBullet.attachToList(listId, 0);\n\n// --- Part 2: delete the first blank spacer paragraph after the PyBank text\nconst pyParas = body.paragraphs;\npyParas.load(\"items/text\");\nawait context.sync();\n\nconst pyScriptsPara = pyParas.items.find(\n  (p) => p.text.indexOf(\"There are two different python scripts here.\") === 0\n);\nif (!pyScriptsPara) {\n  throw new Error(\"Could not find the 'There are two different python scripts...' paragraph\");\n}\n\nconst blankSpacer = pyScriptsPara.getNext();\nblankSpacer.load(\"text\");\nawait context.sync();\nif (blankSpacer.text.trim().length !== 0) {\n  throw new Error(\"Expected the paragraph right after the PyBank text to be blank\");\n}\nblankSpacer.delete();\n\nawait context.sync();\n", "ps1": "# Commit: \"Added cloud computing exp\"\n# 1) Add a new bulleted \"IBM DB2 on Cloud: cloud computing \" line to the\n#    \"Tools:\" list (right after the \"Tableau: data visualizations\" bullet,\n#    continuing the very same numbered list as its sibling bullets).\n# 2) Remove the (now redundant) first blank spacer paragraph that sat right\n#    after the \"There are two different python scripts...\" paragraph.\n\n$doc = $word.ActiveDocument\n\n# --- Part 1: insert the new \"IBM DB2 on Cloud\" bullet ---------------------\n$tableauPara = $null\nforeach ($p in $doc.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"Tableau: data visualizations\") {\n        $tableauPara = $p\n        break\n    }\n}\nif ($tableauPara -eq $null) {\n    throw \"Could not find the 'Tableau: data visualizations' bullet paragraph\"\n}\n\n$listTemplate = $tableauPara.Range.ListFormat.ListTemplate\n\n$tableauPara.Range.InsertParagraphAfter()\n$newPara = $tableauPara.Next()\n$newPara.Range.Text = \"IBM DB2 on Cloud: cloud computing \"\n$newPara.Range.Style = \"List Paragraph\"\n# ContinuePreviousList=$true makes the new bullet join the SAME list (numId)\n# instead of minting a brand-new list instance with identical formatting.\n$newPara.Range.ListFormat.ApplyListTemplate($listTemplate, $true)\n\n# --- Part 2: delete the first blank spacer paragraph after the PyBank text\n$pyScriptsPara = $null\nforeach ($p in $doc.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"There are two different python scripts here.\")) {\n        $pyScriptsPara = $p\n        break\n    }\n}\nif ($pyScriptsPara -eq $null) {\n    throw \"Could not find the 'There are two different python scripts...' paragraph\"\n}\n\n$blankSpacer = $pyScriptsPara.Next()\nif ($blankSpacer.Range.Text.Trim().Length -ne 0) {\n    throw \"Expected the paragraph right after the PyBank text to be blank\"\n}\n$blankSpacer.Range.Delete()\n"}
